# Apply the "spinning mine" tuning edit to mySheet (4) (the first/active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mySheet (4)")
$ws.Activate()

# Row 2 inputs (base values)
$ws.Range("C2").Value = 0.4
$ws.Range("D2").Value = 3

# Row 101 inputs (max values)
$ws.Range("B101").Value = 30
$ws.Range("C101").Value = 0.2
$ws.Range("D101").Value = 0.4
$ws.Range("E101").Value = 15
$ws.Range("F101").Value = 15
$ws.Range("M101").Value = 75

# Move the active cell selection to D12 as recorded in the saved view state.
$ws.Range("D12").Select()
